# Update Mappings 22 Ontologies
# Adds a new "SBO_DEF" column (F) to the BFO/SBO mapping sheet, with an
# empty-list placeholder value "[]" for each of the two data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the new header cell (F1) the same formatting as the other header
# cells (bold font, border, centered) by copying the format from E1 before
# writing the new header text - mirrors the existing B1:E1 style.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("F1").Value = "SBO_DEF"
$ws.Range("F2").Value = "[]"
$ws.Range("F3").Value = "[]"
